$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value for every data row (2..500).
# All of them need to move forward by exactly one day: 45179 -> 45180
# (i.e. 2023-09-10 -> 2023-09-11).
for ($row = 2; $row -le 500; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
